# Applies the cryptos list update (crypto prices / 1h volume changes) for Fri May 31 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to store the exact literal text (avoids Excel auto-converting
    # number-looking strings like "1.00" or "0.520" into floating point numbers,
    # which would silently drop significant trailing/leading characters).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "68.159.47"
Set-TextCell $ws.Range("E2") "  +0.73%  "
Set-TextCell $ws.Range("D3") "3.770.55"
Set-TextCell $ws.Range("E3") "  +1.29%  "
Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  +0.10%  "
Set-TextCell $ws.Range("D5") "593.15"
Set-TextCell $ws.Range("E5") "  +0.25%  "
Set-TextCell $ws.Range("D6") "167.08"
Set-TextCell $ws.Range("E6") "  +1.05%  "
Set-TextCell $ws.Range("D7") "3.769.26"
Set-TextCell $ws.Range("E7") "  +1.26%  "
Set-TextCell $ws.Range("D9") "0.520"
Set-TextCell $ws.Range("E9") "  -0.18%  "
Set-TextCell $ws.Range("D10") "0.159"
Set-TextCell $ws.Range("E10") "  +0.57%  "
Set-TextCell $ws.Range("E11") "  -1.22%  "
Set-TextCell $ws.Range("E12") "  -0.46%  "
Set-TextCell $ws.Range("D13") "0.0000258"
Set-TextCell $ws.Range("E13") "  -1.19%  "
Set-TextCell $ws.Range("D14") "36.03"
Set-TextCell $ws.Range("E14") "  -0.64%  "
Set-TextCell $ws.Range("D15") "4.401.23"
Set-TextCell $ws.Range("E15") "  +1.25%  "
Set-TextCell $ws.Range("D16") "3.756.64"
Set-TextCell $ws.Range("E16") "  +0.84%  "
Set-TextCell $ws.Range("D17") "68.157.83"
Set-TextCell $ws.Range("E17") "  +0.85%  "
Set-TextCell $ws.Range("D18") "17.88"
Set-TextCell $ws.Range("E18") "  -1.94%  "
Set-TextCell $ws.Range("B19") "TRON"
Set-TextCell $ws.Range("C19") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws.Range("D19") "0.112"
Set-TextCell $ws.Range("E19") "  +0.53%  "
Set-TextCell $ws.Range("B20") "Polkadot"
Set-TextCell $ws.Range("C20") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D20") "6.98"
Set-TextCell $ws.Range("E20") "  -0.62%  "
Set-TextCell $ws.Range("D21") "10.75"
Set-TextCell $ws.Range("E21") "  +0.44%  "
Set-TextCell $ws.Range("D22") "464.13"
Set-TextCell $ws.Range("E22") "  -0.43%  "
Set-TextCell $ws.Range("D23") "0.696"
Set-TextCell $ws.Range("E23") "  -0.37%  "
Set-TextCell $ws.Range("D24") "0.0000148"
Set-TextCell $ws.Range("E24") "  +11.40%  "
Set-TextCell $ws.Range("D25") "83.85"
Set-TextCell $ws.Range("E25") "  +1.47%  "
Set-TextCell $ws.Range("E26") "  -1.22%  "
Set-TextCell $ws.Range("D27") "11.82"
Set-TextCell $ws.Range("E27") "  -1.46%  "
Set-TextCell $ws.Range("E28") "  -0.89%  "
Set-TextCell $ws.Range("E29") "  +0.10%  "
Set-TextCell $ws.Range("D30") "2.78"
Set-TextCell $ws.Range("E30") "  +0.40%  "
Set-TextCell $ws.Range("E31") "  -0.30%  "
Set-TextCell $ws.Range("D32") "29.91"
Set-TextCell $ws.Range("E32") "  +1.47%  "
Set-TextCell $ws.Range("D33") "2.14"
Set-TextCell $ws.Range("E33") "  -3.79%  "
Set-TextCell $ws.Range("D34") "9.12"
Set-TextCell $ws.Range("E34") "  +1.22%  "
Set-TextCell $ws.Range("D35") "1.00"
Set-TextCell $ws.Range("E35") "  +0.07%  "
Set-TextCell $ws.Range("D36") "3.725.10"
Set-TextCell $ws.Range("E36") "  +1.38%  "
Set-TextCell $ws.Range("D37") "0.100"
Set-TextCell $ws.Range("E37") "  -0.86%  "
Set-TextCell $ws.Range("D38") "3.45"
Set-TextCell $ws.Range("E38") "  +1.81%  "
Set-TextCell $ws.Range("E39") "  +0.62%  "
Set-TextCell $ws.Range("E40") "  +1.62%  "
Set-TextCell $ws.Range("D41") "5.76"
Set-TextCell $ws.Range("E41") "  +0.34%  "
Set-TextCell $ws.Range("D42") "1.00"
Set-TextCell $ws.Range("E42") "  +0.08%  "
Set-TextCell $ws.Range("D44") "44.42"
Set-TextCell $ws.Range("E44") "  +17.13%  "
Set-TextCell $ws.Range("D45") "0.299"
Set-TextCell $ws.Range("E45") "  -1.75%  "
Set-TextCell $ws.Range("D46") "46.90"
Set-TextCell $ws.Range("E46") "  +3.73%  "
Set-TextCell $ws.Range("D47") "1.90"
Set-TextCell $ws.Range("E47") "  -0.40%  "
Set-TextCell $ws.Range("E48") "  -1.49%  "
Set-TextCell $ws.Range("D49") "145.44"
Set-TextCell $ws.Range("E49") "  +1.59%  "
Set-TextCell $ws.Range("D50") "388.29"
Set-TextCell $ws.Range("E50") "  -0.04%  "
Set-TextCell $ws.Range("D51") "2.785.42"
Set-TextCell $ws.Range("E51") "  +4.25%  "
